$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.189548850059509
$ws.Range("B1").Value = 2.34879732131958
$ws.Range("C1").Value = 3.820623397827148
$ws.Range("D1").Value = 3.101420402526855
$ws.Range("E1").Value = 1.139287352561951
